$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added columns (L and M) to the bead catalog table for rows 58-95,
# each populated with 0.
for ($r = 58; $r -le 95; $r++) {
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
}

# Update the view/selection to reflect the newly added columns.
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 1
$ws.Range("L58:M95").Select()
